$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FABIANA REBECA ARRUNATEGUI SILUPU", 161),
    @("GONZALES FIESTAS MARIA MARIBEL", 152),
    @("ANTON INGA FATIMA DEL ROSARIO", 151),
    @("RUIZ CHIROQUE CLAUDIA JUDITH", 149),
    @("BANCAYAN FIESTA DILVER HUMBERTO", 147),
    @("TEMOCHE ECHE URSULA YESSENIA", 147),
    @("BAUTISTA CHAVESTA ERICKA MEDALIT", 142),
    @("LLENQUE ANTON HELEN JOHANA", 142),
    @("VELASCO PEÑA KAREN ARELLYS", 136),
    @("PINTADO CHASQUERO ESTEFANY", 117),
    @("FLORES SILUPU MARY CARMEN", 116),
    @("ORDINOLA JIBAJA JOSE ALBERTO", 113),
    @("MONDRAGON NONAJULCA MARISOL", 111),
    @("HERNANDEZ CARNERO ARTURO SEBASTIAN", 109),
    @("MORENO YANAYACO NAYLA GUADALUPE", 88),
    @("CASTRO ESTRADA CINTHIA PATRICIA", 84)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
